# This script reproduces a weekly data-refresh commit: two new price
# records (for the week dated 2021-11-10, Excel serial 44508) are inserted
# into the "Hortaliza, Feria Lagunitas de Puerto Montt - Lechuga" table,
# right before the existing row 295. All rows that used to occupy
# 295..373 shift down to 297..375, and the sheet's used-range dimension
# grows from A1:R373 to A1:R375 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two brand-new rows at position 295 (pushes old rows 295-373 down
# to 297-375).
$ws.Range("A295:A296").EntireRow.Insert()

# --- New row 295 -----------------------------------------------------
$ws.Cells.Item(295, 1).Value2 = 4
$ws.Cells.Item(295, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(295, 3).Value2 = "Los Lagos"
$ws.Cells.Item(295, 4).Value2 = 44508
$ws.Cells.Item(295, 5).Value2 = 10
$ws.Cells.Item(295, 6).Value2 = 100112033
$ws.Cells.Item(295, 7).Value2 = "Lechuga"
$ws.Cells.Item(295, 8).Value2 = "Conconina(o)"
$ws.Cells.Item(295, 9).Value2 = "Segunda"
$ws.Cells.Item(295, 10).Value2 = 80
$ws.Cells.Item(295, 11).Value2 = 8000
$ws.Cells.Item(295, 12).Value2 = 8000
$ws.Cells.Item(295, 13).Value2 = 8000
$ws.Cells.Item(295, 14).Value2 = "`$/caja 12 unidades"
$ws.Cells.Item(295, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(295, 16).Value2 = 667
$ws.Cells.Item(295, 17).Value2 = 12
$ws.Cells.Item(295, 18).Value2 = "Hortaliza"

# --- New row 296 -----------------------------------------------------
$ws.Cells.Item(296, 1).Value2 = 4
$ws.Cells.Item(296, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(296, 3).Value2 = "Los Lagos"
$ws.Cells.Item(296, 4).Value2 = 44508
$ws.Cells.Item(296, 5).Value2 = 10
$ws.Cells.Item(296, 6).Value2 = 100112033
$ws.Cells.Item(296, 7).Value2 = "Lechuga"
$ws.Cells.Item(296, 8).Value2 = "Escarola"
$ws.Cells.Item(296, 9).Value2 = "Primera"
$ws.Cells.Item(296, 10).Value2 = 300
$ws.Cells.Item(296, 11).Value2 = 9000
$ws.Cells.Item(296, 12).Value2 = 9000
$ws.Cells.Item(296, 13).Value2 = 9000
$ws.Cells.Item(296, 14).Value2 = "`$/caja 15 unidades"
$ws.Cells.Item(296, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(296, 16).Value2 = 600
$ws.Cells.Item(296, 17).Value2 = 15
$ws.Cells.Item(296, 18).Value2 = "Hortaliza"
